$d = $word.ActiveDocument

# --- Change 1 -------------------------------------------------------------
# "Socks in the Dark" intro paragraph: merge the three runs that spell out
# "... guarantee getting the [following:][ ]" back into a single run and
# drop the gramStart/gramEnd proofErr markers around "following:" by doing
# a same-text Find/Replace that spans the run boundary.
$p1 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -match "guarantee getting the following:") {
        $p1 = $cand
        break
    }
}
$rng1 = $p1.Range
$rng1.Find.Execute("the following: ", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "the following: ", 2) | Out-Null

# --- Change 2 -------------------------------------------------------------
# "Define." paragraph under Part B: merge the three runs around the
# gramStart/gramEnd-wrapped "socks" back into a single run.
$p2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -match "As for obtaining one matching pair") {
        $p2 = $cand
        break
    }
}
$rng2 = $p2.Range
$rng2.Find.Execute("a total of 20 socks in a drawer", $true, $false, $false, `
                    $false, $false, $true, 1, $false, `
                    "a total of 20 socks in a drawer", 2) | Out-Null

# --- Change 3 -------------------------------------------------------------
# "2/5ths of the socks are NOT black." paragraph: turn it into a list item
# (same list/numId as the preceding "Half the socks..." bullet), drop the
# stray paragraph-mark run formatting (sz/szCs 22), and move the _GoBack
# bookmark from the end of the paragraph to its start.
$pPrev = $null
$pTarget = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -match "NOT black") {
        $pTarget = $cand
        $pPrev = $d.Paragraphs.Item($i - 1)
        break
    }
}

# Clear the explicit sz/szCs left on the paragraph mark.
$pTarget.Range.Select()
$word.Selection.ClearFormatting()

# Make it part of the same numbered list as the previous paragraph.
$tmpl = $pPrev.Range.ListFormat.ListTemplate
$pTarget.Range.ListFormat.ApplyListTemplate($tmpl, $true)

# Relocate the _GoBack bookmark to the start of the paragraph.
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()
$startRng = $d.Range($pTarget.Range.Start, $pTarget.Range.Start)
$d.Bookmarks.Add("_GoBack", $startRng)
